# "Update the col name for STH" -- rename the header row of the
# Southbank URL mapping sheet to the lower_snake_case / OSM-style
# column names used by the rest of the mapping pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "building_name"
$ws.Range("B1").Value = "building_no"
$ws.Range("C1").Value = "website:map"

# Column A was re-autofit by Excel after the rename; widen it to fit
# the longest building name in the sheet.
$ws.Columns.Item(1).ColumnWidth = 46.1
